# Auto-generated Excel COM-interop script to apply Goblin_Profits market-data refresh
# Updates cell values across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets to match
# the latest scheduled runner data pull (currentAveragePrice / LevePrice / LeveProfit columns).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 41.916668
$ws.Range("I11").Value = 41.916668
$ws.Range("K11").Value = 41.916668
$ws.Range("M11").Value = 98.083332
$ws.Range("H28").Value = 3928.5862
$ws.Range("I28").Value = 4763.7827
$ws.Range("K28").Value = 4763.7827
$ws.Range("M28").Value = -4278.7827
$ws.Range("H33").Value = 760.3871
$ws.Range("I33").Value = 137.14285
$ws.Range("K33").Value = 137.14285
$ws.Range("M33").Value = 91.85714999999999
$ws.Range("H40").Value = 2470.5881
$ws.Range("J40").Value = 2769.2307
$ws.Range("L40").Value = 2769.2307
$ws.Range("N40").Value = -3119.2307
$ws.Range("H55").Value = 159.3
$ws.Range("I55").Value = 129
$ws.Range("J55").Value = 204.75
$ws.Range("K55").Value = 129
$ws.Range("L55").Value = 204.75
$ws.Range("M55").Value = 85
$ws.Range("N55").Value = -632.75
$ws.Range("H103").Value = 1403.5
$ws.Range("I103").Value = 1347
$ws.Range("J103").Value = 1460
$ws.Range("K103").Value = 4041
$ws.Range("L103").Value = 4380
$ws.Range("M103").Value = -3455
$ws.Range("N103").Value = -5552
$ws.Range("H132").Value = 2387.5557
$ws.Range("J132").Value = 2554.889
$ws.Range("L132").Value = 7664.667
$ws.Range("N132").Value = -12724.667
$ws.Range("H141").Value = 4938.067
$ws.Range("I141").Value = 2257.9443
$ws.Range("J141").Value = 8958.25
$ws.Range("K141").Value = 6773.8329
$ws.Range("L141").Value = 26874.75
$ws.Range("M141").Value = -1593.8329
$ws.Range("N141").Value = -37234.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3565.282
$ws.Range("I32").Value = 3663.5405
$ws.Range("K32").Value = 3663.5405
$ws.Range("M32").Value = -3376.5405
$ws.Range("H132").Value = 2409.3333
$ws.Range("I132").Value = 1704.9286
$ws.Range("K132").Value = 5114.7858
$ws.Range("M132").Value = -2584.7858

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 44763.5
$ws.Range("J26").Value = 78990
$ws.Range("L26").Value = 78990
$ws.Range("N26").Value = -79574
$ws.Range("H99").Value = 3791.96
$ws.Range("I99").Value = 2592.6
$ws.Range("K99").Value = 2592.6
$ws.Range("M99").Value = -1094.6
$ws.Range("H107").Value = 3858.84
$ws.Range("H134").Value = 1768.625
$ws.Range("I134").Value = 1783.5
$ws.Range("J134").Value = 1605
$ws.Range("K134").Value = 5350.5
$ws.Range("L134").Value = 4815
$ws.Range("M134").Value = -2815.5
$ws.Range("N134").Value = -9885

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2245.4
$ws.Range("I2").Value = 2383.2222
$ws.Range("J2").Value = 1005
$ws.Range("K2").Value = 2383.2222
$ws.Range("L2").Value = 1005
$ws.Range("M2").Value = -2270.2222
$ws.Range("N2").Value = -1231
$ws.Range("H16").Value = 2620.4707
$ws.Range("I16").Value = 2544.818
$ws.Range("J16").Value = 2759.1667
$ws.Range("K16").Value = 2544.818
$ws.Range("L16").Value = 2759.1667
$ws.Range("M16").Value = -2257.818
$ws.Range("N16").Value = -3333.1667
$ws.Range("H99").Value = 2420.8572
$ws.Range("I99").Value = 1991
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 1991
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -493
$ws.Range("N99").Value = -7996
$ws.Range("H113").Value = 2620.4707
$ws.Range("I113").Value = 2544.818
$ws.Range("J113").Value = 2759.1667
$ws.Range("K113").Value = 2544.818
$ws.Range("L113").Value = 2759.1667
$ws.Range("M113").Value = -374.8180000000002
$ws.Range("N113").Value = -7099.1667
$ws.Range("H126").Value = 2420.8572
$ws.Range("I126").Value = 1991
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 5973
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3503
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2229.2856
$ws.Range("I132").Value = 2268.25
$ws.Range("J132").Value = 1995.5
$ws.Range("K132").Value = 6804.75
$ws.Range("L132").Value = 5986.5
$ws.Range("M132").Value = -4274.75
$ws.Range("N132").Value = -11046.5
$ws.Range("H134").Value = 2195.125
$ws.Range("J134").Value = 4000
$ws.Range("L134").Value = 12000
$ws.Range("N134").Value = -17070

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 32.48148
$ws.Range("I2").Value = 26.526316
$ws.Range("J2").Value = 46.625
$ws.Range("K2").Value = 159.157896
$ws.Range("L2").Value = 279.75
$ws.Range("M2").Value = -46.15789599999999
$ws.Range("N2").Value = -505.75
$ws.Range("H9").Value = 100140.5
$ws.Range("I9").Value = 75
$ws.Range("J9").Value = 143025.72
$ws.Range("K9").Value = 225
$ws.Range("L9").Value = 429077.16
$ws.Range("M9").Value = -1
$ws.Range("N9").Value = -429525.16
$ws.Range("H86").Value = 266.77777
$ws.Range("I86").Value = 198.33333
$ws.Range("J86").Value = 403.66666
$ws.Range("K86").Value = 594.99999
$ws.Range("L86").Value = 1210.99998
$ws.Range("M86").Value = 591.00001
$ws.Range("N86").Value = -3582.99998
$ws.Range("H89").Value = 266.77777
$ws.Range("I89").Value = 198.33333
$ws.Range("J89").Value = 403.66666
$ws.Range("K89").Value = 1784.99997
$ws.Range("L89").Value = 3632.99994
$ws.Range("M89").Value = 4143.00003
$ws.Range("N89").Value = -15488.99994
$ws.Range("H101").Value = 7999.933
$ws.Range("J101").Value = 7999.933
$ws.Range("L101").Value = 23999.799
$ws.Range("N101").Value = -28867.799
$ws.Range("H120").Value = 51018.332
$ws.Range("J120").Value = 54362.727
$ws.Range("L120").Value = 163088.181
$ws.Range("N120").Value = -172764.181
$ws.Range("H121").Value = 728.7692
$ws.Range("I121").Value = 608.3333
$ws.Range("K121").Value = 1824.9999
$ws.Range("M121").Value = -514.9999
$ws.Range("H140").Value = 2054.5
$ws.Range("J140").Value = 4567
$ws.Range("L140").Value = 13701
$ws.Range("N140").Value = -24061

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 850
$ws.Range("I107").Value = 850
$ws.Range("K107").Value = 850
$ws.Range("M107").Value = 1070
$ws.Range("H113").Value = 40006520
$ws.Range("I113").Value = 100001500
$ws.Range("J113").Value = 9866.666999999999
$ws.Range("K113").Value = 100001500
$ws.Range("L113").Value = 9866.666999999999
$ws.Range("M113").Value = -99999330
$ws.Range("N113").Value = -14206.667
$ws.Range("H122").Value = 7177
$ws.Range("I122").Value = 8003.0967
$ws.Range("K122").Value = 24009.2901
$ws.Range("M122").Value = -21559.2901
$ws.Range("H126").Value = 2706.25
$ws.Range("I126").Value = 2706.25
$ws.Range("K126").Value = 8118.75
$ws.Range("M126").Value = -5648.75
$ws.Range("H132").Value = 4518.75
$ws.Range("I132").Value = 4450
$ws.Range("K132").Value = 13350
$ws.Range("M132").Value = -10820

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3097.3
$ws.Range("I7").Value = 3252.5557
$ws.Range("J7").Value = 1700
$ws.Range("K7").Value = 3252.5557
$ws.Range("L7").Value = 1700
$ws.Range("M7").Value = -3140.5557
$ws.Range("N7").Value = -1924
$ws.Range("H46").Value = 966.6667
$ws.Range("I46").Value = 899
$ws.Range("J46").Value = 1000.5
$ws.Range("K46").Value = 899
$ws.Range("L46").Value = 1000.5
$ws.Range("M46").Value = -711
$ws.Range("N46").Value = -1376.5
$ws.Range("H61").Value = 4101.227
$ws.Range("I61").Value = 1665
$ws.Range("K61").Value = 1665
$ws.Range("M61").Value = -1463
$ws.Range("H113").Value = 4101.227
$ws.Range("I113").Value = 1665
$ws.Range("K113").Value = 1665
$ws.Range("M113").Value = 505
$ws.Range("H119").Value = 67500
$ws.Range("J119").Value = 67500
$ws.Range("L119").Value = 67500
$ws.Range("N119").Value = -77176
$ws.Range("H126").Value = 3097.3
$ws.Range("I126").Value = 3252.5557
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 9757.667099999999
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -7287.667099999999
$ws.Range("N126").Value = -10040

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 79999
$ws.Range("J27").Value = 79999
$ws.Range("L27").Value = 79999
$ws.Range("N27").Value = -80137
$ws.Range("H63").Value = 21333
$ws.Range("J63").Value = 21333
$ws.Range("L63").Value = 21333
$ws.Range("N63").Value = -22581
$ws.Range("H66").Value = 21333
$ws.Range("J66").Value = 21333
$ws.Range("L66").Value = 63999
$ws.Range("N66").Value = -70239
$ws.Range("H100").Value = 807.8
$ws.Range("I100").Value = 807.8
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1615.6
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1074.6
$ws.Range("N100").Value = $null
$ws.Range("H113").Value = 1500.6818
$ws.Range("I113").Value = 1414.0625
$ws.Range("J113").Value = 1731.6666
$ws.Range("K113").Value = 4242.1875
$ws.Range("L113").Value = 5194.9998
$ws.Range("M113").Value = -2072.1875
$ws.Range("N113").Value = -9534.9998
$ws.Range("H115").Value = 94982.336
$ws.Range("J115").Value = 94982.336
$ws.Range("L115").Value = 94982.336
$ws.Range("N115").Value = -98116.336
$ws.Range("H122").Value = 4596.3887
$ws.Range("I122").Value = 973.5
$ws.Range("K122").Value = 2920.5
$ws.Range("M122").Value = -470.5
$ws.Range("H132").Value = 5041.533
$ws.Range("I132").Value = 5355.615
$ws.Range("K132").Value = 16066.845
$ws.Range("M132").Value = -13536.845
